# Macro Schedule/Details import: add "District ID" column (F) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header cell F1 with shared string "District ID"
$ws.Range("F1").Value = "District ID"

# Column F was auto best-fit sized in the original edit to fit the new header text.
# ColumnWidth is specified in characters; closest achievable value to the
# authored stored width (9.5703125) given this runtime's rounding.
$ws.Columns.Item(6).ColumnWidth = 8.74

# The author's last on-screen action was clicking the column G header,
# selecting the (empty) next column.
$ws.Columns.Item(7).Select() | Out-Null
